$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 40
$ws_ALC.Range("H40").Value = 3856.5293
$ws_ALC.Range("I40").Value = 3044.7317
$ws_ALC.Range("J40").Value = 7184.9
$ws_ALC.Range("K40").Value = 3044.7317
$ws_ALC.Range("L40").Value = 7184.9
$ws_ALC.Range("M40").Value = -2869.7317
$ws_ALC.Range("N40").Value = -7534.9

# ALC row 46
$ws_ALC.Range("H46").Value = 2000
$ws_ALC.Range("J46").Value = 2000
$ws_ALC.Range("L46").Value = 6000
$ws_ALC.Range("N46").Value = -6238

# ALC row 60
$ws_ALC.Range("H60").Value = 2000
$ws_ALC.Range("J60").Value = 2000
$ws_ALC.Range("L60").Value = 6000
$ws_ALC.Range("N60").Value = -6968

# ALC row 69
$ws_ALC.Range("H69").Value = 6772.375
$ws_ALC.Range("J69").Value = 7096.0713
$ws_ALC.Range("L69").Value = 21288.2139
$ws_ALC.Range("N69").Value = -23036.2139

# ALC row 72
$ws_ALC.Range("H72").Value = 6772.375
$ws_ALC.Range("J72").Value = 7096.0713
$ws_ALC.Range("L72").Value = 63864.64169999999
$ws_ALC.Range("N72").Value = -72600.64169999999

# ALC row 80
$ws_ALC.Range("H80").Value = 2314.7144
$ws_ALC.Range("J80").Value = 1666.3334
$ws_ALC.Range("L80").Value = 4999.0002
$ws_ALC.Range("N80").Value = -6995.0002

# ALC row 83
$ws_ALC.Range("H83").Value = 2314.7144
$ws_ALC.Range("J83").Value = 1666.3334
$ws_ALC.Range("L83").Value = 14997.0006
$ws_ALC.Range("N83").Value = -24981.0006

# ALC row 132
$ws_ALC.Range("H132").Value = 15054.4375
$ws_ALC.Range("I132").Value = 15391
$ws_ALC.Range("J132").Value = 10006
$ws_ALC.Range("K132").Value = 46173
$ws_ALC.Range("L132").Value = 30018
$ws_ALC.Range("M132").Value = -43643
$ws_ALC.Range("N132").Value = -35078

# ALC row 135
$ws_ALC.Range("H135").Value = 982.1667
$ws_ALC.Range("I135").Value = 744.2727
$ws_ALC.Range("J135").Value = 3599
$ws_ALC.Range("K135").Value = 6698.454299999999
$ws_ALC.Range("L135").Value = 32391
$ws_ALC.Range("M135").Value = -4163.454299999999
$ws_ALC.Range("N135").Value = -37461

# ALC row 137
$ws_ALC.Range("H137").Value = 2960.7708
$ws_ALC.Range("J137").Value = 3763.68
$ws_ALC.Range("L137").Value = 11291.04
$ws_ALC.Range("N137").Value = -16391.04

# ALC row 138
$ws_ALC.Range("H138").Value = 2317.7334
$ws_ALC.Range("I138").Value = 1329.4546
$ws_ALC.Range("J138").Value = 5035.5
$ws_ALC.Range("K138").Value = 3988.3638
$ws_ALC.Range("L138").Value = 15106.5
$ws_ALC.Range("M138").Value = 1151.6362
$ws_ALC.Range("N138").Value = -25386.5

# ARM row 61
$ws_ARM.Range("H61").Value = 3873.9167
$ws_ARM.Range("I61").Value = 3148.7
$ws_ARM.Range("K61").Value = 3148.7
$ws_ARM.Range("M61").Value = -2936.7

# ARM row 132
$ws_ARM.Range("H132").Value = 5782.2354
$ws_ARM.Range("I132").Value = 5307.0713
$ws_ARM.Range("J132").Value = 7999.6665
$ws_ARM.Range("K132").Value = 15921.2139
$ws_ARM.Range("L132").Value = 23998.9995
$ws_ARM.Range("M132").Value = -13391.2139
$ws_ARM.Range("N132").Value = -29058.9995

# ARM row 136
$ws_ARM.Range("H136").Value = 3873.9167
$ws_ARM.Range("I136").Value = 3148.7
$ws_ARM.Range("K136").Value = 9446.099999999999
$ws_ARM.Range("M136").Value = -6896.099999999999

# BSM row 100
$ws_BSM.Range("H100").Value = 5320
$ws_BSM.Range("J100").Value = 5320
$ws_BSM.Range("L100").Value = 5320
$ws_BSM.Range("N100").Value = -7484

# CRP row 31
$ws_CRP.Range("H31").Value = 4656.193
$ws_CRP.Range("J31").Value = 5699.9556
$ws_CRP.Range("L31").Value = 5699.9556
$ws_CRP.Range("N31").Value = -6289.9556

# CRP row 34
$ws_CRP.Range("H34").Value = 4656.193
$ws_CRP.Range("J34").Value = 5699.9556
$ws_CRP.Range("L34").Value = 5699.9556
$ws_CRP.Range("N34").Value = -6103.9556

# CRP row 68
$ws_CRP.Range("H68").Value = 64382.5
$ws_CRP.Range("J68").Value = 64382.5
$ws_CRP.Range("L68").Value = 64382.5
$ws_CRP.Range("N68").Value = -65880.5

# CRP row 71
$ws_CRP.Range("H71").Value = 64382.5
$ws_CRP.Range("J71").Value = 64382.5
$ws_CRP.Range("L71").Value = 193147.5
$ws_CRP.Range("N71").Value = -200635.5

# CRP row 99
$ws_CRP.Range("H99").Value = 2400
$ws_CRP.Range("I99").Value = 0
$ws_CRP.Range("K99").Value = 0
$ws_CRP.Range("M99").ClearContents()

# CRP row 106
$ws_CRP.Range("H106").Value = 27667.5
$ws_CRP.Range("J106").Value = 27667.5
$ws_CRP.Range("L106").Value = 27667.5
$ws_CRP.Range("N106").Value = -30191.5

# CRP row 126
$ws_CRP.Range("H126").Value = 2400
$ws_CRP.Range("I126").Value = 0
$ws_CRP.Range("K126").Value = 0
$ws_CRP.Range("M126").ClearContents()

# CUL row 52
$ws_CUL.Range("H52").Value = 532
$ws_CUL.Range("J52").Value = 532
$ws_CUL.Range("L52").Value = 1596
$ws_CUL.Range("N52").Value = -2128

# GSM row 63
$ws_GSM.Range("H63").Value = 48330
$ws_GSM.Range("J63").Value = 50000
$ws_GSM.Range("L63").Value = 50000
$ws_GSM.Range("N63").Value = -51372

# GSM row 66
$ws_GSM.Range("H66").Value = 48330
$ws_GSM.Range("J66").Value = 50000
$ws_GSM.Range("L66").Value = 150000
$ws_GSM.Range("N66").Value = -156864

# GSM row 80
$ws_GSM.Range("H80").Value = 3111.625
$ws_GSM.Range("I80").Value = 2641
$ws_GSM.Range("K80").Value = 2641
$ws_GSM.Range("M80").Value = -1643

# GSM row 83
$ws_GSM.Range("H83").Value = 3111.625
$ws_GSM.Range("I83").Value = 2641
$ws_GSM.Range("K83").Value = 13205
$ws_GSM.Range("M83").Value = -8213

# GSM row 105
$ws_GSM.Range("H105").Value = 15924.833
$ws_GSM.Range("J105").Value = 15924.833
$ws_GSM.Range("L105").Value = 15924.833
$ws_GSM.Range("N105").Value = -22912.833

# GSM row 107
$ws_GSM.Range("H107").Value = 1062.8
$ws_GSM.Range("I107").Value = 1121
$ws_GSM.Range("J107").Value = 902.75
$ws_GSM.Range("K107").Value = 1121
$ws_GSM.Range("L107").Value = 902.75
$ws_GSM.Range("M107").Value = 799
$ws_GSM.Range("N107").Value = -4742.75

# GSM row 122
$ws_GSM.Range("H122").Value = 3333
$ws_GSM.Range("I122").Value = 3333
$ws_GSM.Range("K122").Value = 9999
$ws_GSM.Range("M122").Value = -7549

# GSM row 132
$ws_GSM.Range("H132").Value = 2968.5334
$ws_GSM.Range("I132").Value = 2537.7856
$ws_GSM.Range("J132").Value = 8999
$ws_GSM.Range("K132").Value = 7613.3568
$ws_GSM.Range("L132").Value = 26997
$ws_GSM.Range("M132").Value = -5083.3568
$ws_GSM.Range("N132").Value = -32057

# LTW row 22
$ws_LTW.Range("H22").Value = 1658.8334
$ws_LTW.Range("I22").Value = 1481.8572
$ws_LTW.Range("J22").Value = 2278.25
$ws_LTW.Range("K22").Value = 1481.8572
$ws_LTW.Range("L22").Value = 2278.25
$ws_LTW.Range("M22").Value = -1186.8572
$ws_LTW.Range("N22").Value = -2868.25

# LTW row 27
$ws_LTW.Range("H27").Value = 1658.8334
$ws_LTW.Range("I27").Value = 1481.8572
$ws_LTW.Range("J27").Value = 2278.25
$ws_LTW.Range("K27").Value = 1481.8572
$ws_LTW.Range("L27").Value = 2278.25
$ws_LTW.Range("M27").Value = -1374.8572
$ws_LTW.Range("N27").Value = -2492.25

# LTW row 46
$ws_LTW.Range("H46").Value = 5594
$ws_LTW.Range("I46").Value = 4830.1665
$ws_LTW.Range("J46").Value = 6739.75
$ws_LTW.Range("K46").Value = 4830.1665
$ws_LTW.Range("L46").Value = 6739.75
$ws_LTW.Range("M46").Value = -4642.1665
$ws_LTW.Range("N46").Value = -7115.75

# LTW row 100
$ws_LTW.Range("H100").Value = 7627.1816
$ws_LTW.Range("I100").Value = 4779.8
$ws_LTW.Range("K100").Value = 4779.8
$ws_LTW.Range("M100").Value = -4238.8

# LTW row 104
$ws_LTW.Range("H104").Value = 8950
$ws_LTW.Range("J104").Value = 8950
$ws_LTW.Range("L104").Value = 8950
$ws_LTW.Range("N104").Value = -15938

# LTW row 136
$ws_LTW.Range("H136").Value = 2574.2
$ws_LTW.Range("I136").Value = 1706
$ws_LTW.Range("J136").Value = 4600
$ws_LTW.Range("K136").Value = 5118
$ws_LTW.Range("L136").Value = 13800
$ws_LTW.Range("M136").Value = -2568
$ws_LTW.Range("N136").Value = -18900

# LTW row 137
$ws_LTW.Range("H137").Value = 429999
$ws_LTW.Range("J137").Value = 429999
$ws_LTW.Range("L137").Value = 429999
$ws_LTW.Range("N137").Value = -440199

# WVR row 122
$ws_WVR.Range("H122").Value = 5180
$ws_WVR.Range("I122").Value = 5000
$ws_WVR.Range("J122").Value = 5225
$ws_WVR.Range("K122").Value = 15000
$ws_WVR.Range("L122").Value = 15675
$ws_WVR.Range("M122").Value = -12550
$ws_WVR.Range("N122").Value = -20575
